$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Columns whose values are swapped between row 12 and row 14
$cols = @("A", "B", "E", "F", "G", "H", "P", "Q", "R", "S")

foreach ($col in $cols) {
    $range12 = $ws.Range($col + "12")
    $range14 = $ws.Range($col + "14")

    $v12 = $range12.Value2
    $v14 = $range14.Value2

    $range12.Value = $v14
    $range14.Value = $v12
}
